$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.859.21'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.79%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.766.23'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +7.08%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '421.32'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.34'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.756.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +7.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.653'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.779'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.188'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +15.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000435'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +62.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.96'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.43'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.355.20'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +7.01%  '
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.753.32'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +6.20%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.66'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.10'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.14'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.746.53'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '452.56'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.92'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +18.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.44'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.27%  '
$ws.Range("E25").Value = '  -4.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '39.08'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +13.72%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.33'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.45%  '
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.16'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.10'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.60%  '
$ws.Range("E30").Value = '  +6.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.75'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.77'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.25'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.65%  '
$ws.Range("E34").Value = '  +2.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '42.05'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +5.07%  '
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  -3.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0770'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.44%  '
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.98'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +27.07%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.996'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '28.02'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +28.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.42'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.60%  '
$ws.Range("E45").Value = '  +6.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '148.10'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.18'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +23.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.92'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.36%  '
$ws.Range("E49").Value = '  -3.87%  '
$ws.Range("E50").Value = '  -4.82%  '
$ws.Range("E51").Value = '  -1.89%  '
